$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.628.00"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -5.34%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.010.98"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -6.85%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "579.32"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "126.88"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -8.26%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.007.44"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -6.90%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.500"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.21%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -7.94%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.16"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.64%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.441"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.03%  "
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -7.94%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.67"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -8.23%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.20%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.502.19"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -7.00%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.006.80"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -6.85%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "60.538.51"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.54%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -4.09%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "432.61"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -7.74%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.13"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -7.33%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.666"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -6.37%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -9.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.80"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -5.59%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "79.60"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -5.23%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.998"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -0.19%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -4.97%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -7.86%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -8.79%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.19"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -10.54%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -8.84%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0935"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -10.25%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -13.29%  "
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -8.62%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -5.89%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "50.04"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.33%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -10.57%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "8.39"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +2.35%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0361"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -9.15%  "
$ws.Range("B41").NumberFormat = "@"
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").NumberFormat = "@"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "387.23"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -5.87%  "
$ws.Range("B42").NumberFormat = "@"
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").NumberFormat = "@"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.110"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -3.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.50"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -10.81%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.660.20"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -7.18%  "
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.20%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -8.69%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.02"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -7.95%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "118.40"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -7.83%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -4.60%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "23.77"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -8.68%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.135"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +2.38%  "
